$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.991.53"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "1.883.29"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "312.86"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  -3.67%  "
$ws.Range("D8").Value = "0.3848"
$ws.Range("E8").Value = "  -2.89%  "
$ws.Range("D9").Value = "0.09190"
$ws.Range("E9").Value = "  -5.69%  "
$ws.Range("E10").Value = "  -2.92%  "
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "6.330"
$ws.Range("E12").Value = "  -3.41%  "
$ws.Range("D13").Value = "20.73"
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("D14").Value = "1.875.56"
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("D15").Value = "7.284"
$ws.Range("E15").Value = "  -3.32%  "
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "0.00001105"
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("D18").Value = "91.29"
$ws.Range("E18").Value = "  -3.76%  "
$ws.Range("D19").Value = "0.06632"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "17.97"
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "6.168"
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("D23").Value = "28.041.83"
$ws.Range("E23").Value = "  -2.25%  "
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("D25").Value = "2.306"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").Value = "2.092.95"
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("D27").Value = "2.540"
$ws.Range("E27").Value = "  -5.35%  "
$ws.Range("D28").Value = "157.75"
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D29").Value = "20.76"
$ws.Range("E29").Value = "  -2.41%  "
$ws.Range("D30").Value = "126.56"
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("D31").Value = "1.067"
$ws.Range("E31").Value = "  -4.60%  "
$ws.Range("D32").Value = "0.1054"
$ws.Range("E32").Value = "  -2.71%  "
$ws.Range("D33").Value = "5.584"
$ws.Range("E33").Value = "  -3.32%  "
$ws.Range("E34").Value = "  -1.22%  "
$ws.Range("D35").Value = "9.365"
$ws.Range("E35").Value = "  -5.70%  "
$ws.Range("D36").Value = "0.06570"
$ws.Range("E36").Value = "  -3.36%  "
$ws.Range("D37").Value = "0.02399"
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("D39").Value = "1.286"
$ws.Range("E39").Value = "  +7.86%  "
$ws.Range("E40").Value = "  -4.73%  "
$ws.Range("E41").Value = "  -2.07%  "
$ws.Range("D42").Value = "0.6409"
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("D43").Value = "4.933"
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("D45").Value = "13.29"
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("D46").Value = "0.6019"
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").Value = "1.293"
$ws.Range("E47").Value = "  +0.89%  "
$ws.Range("D48").Value = "3.672"
$ws.Range("E48").Value = "  -2.44%  "
$ws.Range("D49").Value = "1.989"
$ws.Range("E49").Value = "  -2.45%  "
$ws.Range("D50").Value = "1.200"
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("D51").Value = "121.27"
$ws.Range("E51").Value = "  -3.17%  "
